# Raitha Dinachari.xlsx — "Add files via upload" commit reproduction
#
# Summary of the edit:
#  1. "Daily Expenditure" sheet: two pending blank placeholder rows (old
#     rows 72 & 73, both "65 / 31-May-2022 / Income / (blank)") get filled
#     in with real transactions (Ginger/Medicine expense, and a big Salary
#     income), and two *new* blank placeholder rows are appended after them
#     (new rows 74 & 75) ready for future entries — row 75 is an exact copy
#     of the pattern the old rows 72/73 used.
#  2. A brand new "Ganapati" worksheet is added at the end of the workbook
#     with a small collection-tracking table (Name list, blank Amount /
#     Transaction columns to be filled in later).

# Switch to manual calculation while we make all the edits — this workbook
# has a lot of volatile (TODAY()-based) interest-calculation formulas and
# recalculating after every single write is extremely slow. We force one
# full recalculation right before saving instead.
$excel.Calculation = -4135

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Daily Expenditure" — fill in rows 72 & 73, append rows 74 & 75
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("Daily Expenditure")

# Row 72 was a blank placeholder (65 / 31-May-2022 / Income / ---).
# Turn it into a real expense entry: Ginger / Medicine, 5000.
$de.Range("A72").Value = 66
$de.Range("B72").Value = 44763          # 21-Jul-2022
$de.Range("E72").Value = "Expense"
$de.Range("F72").Value = "Ginger"
$de.Range("G72").Value = "Medicine"
$de.Range("H72").Value = 5000

# Row 73 was also a blank placeholder; turn it into a Salary income entry.
$de.Range("B73").Value = 44770          # 28-Jul-2022
$de.Range("F73").Value = "Salary"
$de.Range("F73").Font.Bold = $true
$de.Range("G73").Value = "July Salary"
$de.Range("H73").Value = 230601

# New row 74: a fresh blank "Expense" placeholder dated 29-Jul-2022.
$de.Range("A74").Value = 65
$de.Range("B74").Value = 44771          # 29-Jul-2022
$de.Range("C74").Formula = '=TEXT(B74,"mmm")'
$de.Range("D74").Formula = '=TEXT(B74,"yyyy")'
$de.Range("E74").Value = "Expense"

# New row 75: re-creates the old blank "Income" placeholder pattern
# (65 / 31-May-2022 / Income / ---) that rows 72/73 used to hold.
$de.Range("A75").Value = 65
$de.Range("B75").Value = 44712          # 31-May-2022
$de.Range("C75").Formula = '=TEXT(B75,"mmm")'
$de.Range("D75").Formula = '=TEXT(B75,"yyyy")'
$de.Range("E75").Value = "Income"

# ---------------------------------------------------------------------
# 2. Add the new "Ganapati" worksheet at the end of the workbook
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$gp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$gp.Name = "Ganapati"

$gp.Range("A1").Value = "SL No"
$gp.Range("B1").Value = "Name"
$gp.Range("C1").Value = "Amount "
$gp.Range("D1").Value = "Transaction"

$gp.Range("A3").Value = 1
$gp.Range("B3").Value = "Koushik"
$gp.Range("A4").Value = 2
$gp.Range("B4").Value = "Dhayanand"
$gp.Range("A5").Value = 3
$gp.Range("B5").Value = "Nithin"
$gp.Range("A6").Value = 4
$gp.Range("B6").Value = "Akshay"
$gp.Range("A7").Value = 5
$gp.Range("B7").Value = "Harish"
$gp.Range("A8").Value = 6
$gp.Range("B8").Value = "Chiru"
$gp.Range("A9").Value = 7
$gp.Range("B9").Value = "Shashank"

$gp.Columns.Item(1).ColumnWidth = 4.9765625
$gp.Columns.Item(2).ColumnWidth = 10.0859375
$gp.Columns.Item(4).ColumnWidth = 10.35546875

# ---------------------------------------------------------------------
# 3. Restore "Daily Expenditure" as the active sheet/selection
# ---------------------------------------------------------------------
$de.Activate()
$de.Range("F74").Select()

# Recalculate once, now that every edit has been made, then return to
# automatic mode (matches the workbook's original calcPr settings).
$excel.Calculate()
$excel.Calculation = -4105
